# Generate Report for Handback
# Appends a new handback row (02f61d55-cf2a-49f8-bf3b-12478da1fb15) to the
# Overview / zh-cn / de-de report tables, mirroring the existing
# b9c3f894-3176-450f-8722-bfc3542103f1 row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet (new row 3)
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ovTable = $ov.ListObjects.Item(1)
$ovTable.ListRows.Add() | Out-Null

$ov.Range("A3").Value = "02f61d55-cf2a-49f8-bf3b-12478da1fb15.md"
$ov.Range("B3").Value = "e2e\02f61d55-cf2a-49f8-bf3b-12478da1fb15.md"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = "2016-09-09 10:08:58"

$ov.Range("B3").Style = "HyperLink"
$ov.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad2a512109d5bd82735f57bd8706c3a6525ce8d7/e2e/02f61d55-cf2a-49f8-bf3b-12478da1fb15.md", "", "", "e2e\02f61d55-cf2a-49f8-bf3b-12478da1fb15.md") | Out-Null

# ---------------------------------------------------------------------
# zh-cn sheet (new row 3)
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zhTable = $zh.ListObjects.Item(1)
$zhTable.ListRows.Add() | Out-Null

$zh.Range("A3").Value = "02f61d55-cf2a-49f8-bf3b-12478da1fb15.md"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "'True"
$zh.Range("G3").Value = "02f61d55-cf2a-49f8-bf3b-12478da1fb15.fbfd7004a46e47e269c995ab04aaa2904503373f.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-09 10:08:46"
$zh.Range("I3").Value = "02f61d55-cf2a-49f8-bf3b-12478da1fb15.md"
$zh.Range("J3").Value = "02f61d55-cf2a-49f8-bf3b-12478da1fb15.fbfd7004a46e47e269c995ab04aaa2904503373f.zh-cn.xlf"
$zh.Range("K3").Value = "2016-09-09 10:09:37"
$zh.Range("L3").Value = "'"
$zh.Range("M3").Value = "'True"
$zh.Range("N3").Value = "'"
$zh.Range("O3").Value = "'False"
$zh.Range("P3").Value = "'"

$zh.Range("A3").Style = "HyperLink"
$zh.Range("I3").Style = "HyperLink"
$zh.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad2a512109d5bd82735f57bd8706c3a6525ce8d7/e2e/02f61d55-cf2a-49f8-bf3b-12478da1fb15.md", "", "", "02f61d55-cf2a-49f8-bf3b-12478da1fb15.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/91456888335e6e2b332c3390eda343a37bccca00/e2e/02f61d55-cf2a-49f8-bf3b-12478da1fb15.md", "", "", "02f61d55-cf2a-49f8-bf3b-12478da1fb15.md") | Out-Null

# ---------------------------------------------------------------------
# de-de sheet (new row 3)
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$deTable = $de.ListObjects.Item(1)
$deTable.ListRows.Add() | Out-Null

$de.Range("A3").Value = "02f61d55-cf2a-49f8-bf3b-12478da1fb15.md"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "'True"
$de.Range("G3").Value = "02f61d55-cf2a-49f8-bf3b-12478da1fb15.fbfd7004a46e47e269c995ab04aaa2904503373f.de-de.xlf"
$de.Range("H3").Value = "2016-09-09 10:08:58"
$de.Range("I3").Value = "02f61d55-cf2a-49f8-bf3b-12478da1fb15.md"
$de.Range("J3").Value = "02f61d55-cf2a-49f8-bf3b-12478da1fb15.fbfd7004a46e47e269c995ab04aaa2904503373f.de-de.xlf"
$de.Range("K3").Value = "2016-09-09 10:09:55"
$de.Range("L3").Value = "'"
$de.Range("M3").Value = "'True"
$de.Range("N3").Value = "'"
$de.Range("O3").Value = "'False"
$de.Range("P3").Value = "'"

$de.Range("A3").Style = "HyperLink"
$de.Range("I3").Style = "HyperLink"
$de.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ad2a512109d5bd82735f57bd8706c3a6525ce8d7/e2e/02f61d55-cf2a-49f8-bf3b-12478da1fb15.md", "", "", "02f61d55-cf2a-49f8-bf3b-12478da1fb15.md") | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/aa9ddb5f27d56c1baf2d7a730fbc74483749a496/e2e/02f61d55-cf2a-49f8-bf3b-12478da1fb15.md", "", "", "02f61d55-cf2a-49f8-bf3b-12478da1fb15.md") | Out-Null

$ov.Select()
$ov.Range("A1").Select()
